# 1) Updated Test Data File / Updated Test for sentiment category:
#    populate the Sentiment_Filters lookup sheet with the new
#    "Sentiment Category" values (Overall/Atmosphere/Product/
#    Customer Service/Value) in column B.
$wb = $excel.ActiveWorkbook

$wsSentiment = $wb.Worksheets.Item("Sentiment_Filters")
$wsSentiment.Range("B2").Value = "Overall"
$wsSentiment.Range("B3").Value = "Atmosphere"
$wsSentiment.Range("B4").Value = "Product"
$wsSentiment.Range("B5").Value = "Customer Service"
$wsSentiment.Range("B6").Value = "Value"
$wsSentiment.Range("B6").Select() | Out-Null

# 2) Updated Sentiment Category code to verify:
#    add a new "MultiSentiment" column (F) to the Reviews_AdvancedFilters
#    test data and flip the join value in column E from "null" to "and" so
#    the advanced-filter scenario exercises the new sentiment-category
#    (multi-select) logic.
$wsAdv = $wb.Worksheets.Item("Reviews_AdvancedFilters")
$wsAdv.Activate() | Out-Null

$wsAdv.Range("F1").Value = "MultiSentiment"
$wsAdv.Range("F1").Font.Bold = $true
$wsAdv.Range("E2").Value = "and"
$wsAdv.Range("F2").Value = "Overall,Atmosphere"

$wsAdv.Columns.Item(6).ColumnWidth = 14.5

$wsAdv.Range("F2").Select() | Out-Null
